$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Location 2: paragraph "Crear CU Ejecutar Spuffy" -> prefix with "[OK] "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Crear CU Ejecutar", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[OK] Crear CU Ejecutar", 2)

$pCrear = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "[OK] Crear CU Ejecutar*") {
        $pCrear = $d.Paragraphs.Item($i)
        break
    }
}
$crearStart = $pCrear.Range.Start

# Split "[OK] " away from "Crear CU Ejecutar " (same rPr, just a run boundary)
# by dropping a bookmark at the seam and removing it again.
$splitCrear = $d.Range($crearStart + 5, $crearStart + 5)
$d.Bookmarks.Add("zzzTmpSplitCrear", $splitCrear)
$d.Bookmarks.Item("zzzTmpSplitCrear").Delete()

# ---------------------------------------------------------------------------
# Location 1: paragraph "Gestión de excepciones: ..." -> prefix with "[OK] "
# and the caret (now at "_GoBack") ends up between "[OK" and "] "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Gestión de excepciones", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[OK] Gestión de excepciones", 2)

$pGestion = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "[OK] Gesti*n de excepciones*") {
        $pGestion = $d.Paragraphs.Item($i)
        break
    }
}
$gestionStart = $pGestion.Range.Start

# Split off "] " from "Gestión de excepciones..." first (temp bookmark, removed).
$splitGestionB = $d.Range($gestionStart + 5, $gestionStart + 5)
$d.Bookmarks.Add("zzzTmpSplitGestion", $splitGestionB)
$d.Bookmarks.Item("zzzTmpSplitGestion").Delete()

# Now split "[OK" from "] " and drop the (moved) _GoBack bookmark exactly there.
# Re-adding "_GoBack" here also removes it from its old location at doc end.
$splitGestionA = $d.Range($gestionStart + 3, $gestionStart + 3)
$d.Bookmarks.Add("_GoBack", $splitGestionA)

Write-Output "Done"
